$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(4)
$p.Range.InsertParagraphBefore()
$p1 = $d.Paragraphs.Item(4)
$p1.Range.Text = "Sky must fade to fog color at horizon"

$p2target = $d.Paragraphs.Item(5)
$p2target.Range.InsertParagraphBefore()
$p2 = $d.Paragraphs.Item(5)
$p2.Range.Text = "Water shader needs to support fog as well"

Write-Output "done"
